# Sprint109 manual testcases update:
# Insert a new row (new row 19) between the existing "Click the Filter->Export" (row18)
# and "Click the Filter->Suppliers" (old row19, becomes row20) testcase rows, describing
# the "Filter->Outlets" testcase, and nudge the sheet's active selection/view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at position 19 (pushes old rows 19-22 down to 20-23) ---
$ws.Rows("19:19").Insert()

# The freshly inserted row's A/D/E/F cells already inherit the correct column
# styles (13/10/10/6) from the row below, but B/C/G lose their border/font
# formatting on insert - restore it by copying format from the row right below
# (which still carries the original, untouched formatting).
$ws.Range("B20:C20").Copy()
$ws.Range("B19:C19").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G20").Copy()
$ws.Range("G19").PasteSpecial(-4122)       # xlPasteFormats

# Row height for the new row
$ws.Rows("19:19").RowHeight = 57.75

# --- 2. Fill in the new row's values ---
# Column A (SL. No) is intentionally left blank, matching the target state.
$ws.Range("B19").Value = "Buyers Home page"
$ws.Range("C19").Value = "Reports menu page"

# Column D: "Click the " + bold("Filter->Outlets")
$ws.Range("D19").Value = "Click the Filter->Outlets"
$d19Bold = $ws.Range("D19").Characters(11, 16)
$d19Bold.Font.Name = "Calibri"
$d19Bold.Font.Size = 11
$d19Bold.Font.Color = 0
$d19Bold.Font.Bold = $true

# Column E: "It should show List of Outlets and" + bold(" ") + "'" + bold("Select all and Deselect all and Apply or Reset'")
$ws.Range("E19").Value = "It should show List of Outlets and 'Select all and Deselect all and Apply or Reset'"
$e19Bold1 = $ws.Range("E19").Characters(35, 1)
$e19Bold1.Font.Name = "Calibri"
$e19Bold1.Font.Size = 11
$e19Bold1.Font.Color = 0
$e19Bold1.Font.Bold = $true
$e19Bold2 = $ws.Range("E19").Characters(37, 47)
$e19Bold2.Font.Name = "Calibri"
$e19Bold2.Font.Size = 11
$e19Bold2.Font.Color = 0
$e19Bold2.Font.Bold = $true

# Column F: plain expected-output text
$ws.Range("F19").Value = "1.Once click the selected option and Apply It will display the Details of selected outlets                                                                             2.if you select Reset it will back to normal page"

# Column G: Result = Pass
$ws.Range("G19").Value = "Pass"

# --- 3. Update the sheet view: selection moved to G20, scrolled down a bit ---
$ws.Range("G20").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 4
